# Update cryptos list figures (price and 1h volume change) per Aug 10 2023 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.476.22"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.850.21"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.80"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6279"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.48%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07734"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.862.75"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6907"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.000"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.99%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009754"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.113.71"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.229"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.531.70"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "233.01"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.28%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9999"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.623"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.50"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1386"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.454"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.475"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05919"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.08%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.106"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.026"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.875"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.168"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7183"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.591"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.84%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.794"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.238.73"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01791"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9083"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.119"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.28%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.027.64"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.32"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "67.18"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.379"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +9.70%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4038"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.128"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.71%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.13%  "
